# Applies the crypto price/volume updates described by the commit diff.
# Numeric-looking text values are prefixed with a literal apostrophe so Excel
# keeps storing them as text (matching the original inlineStr cell type)
# instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '30.569.62'),
    @('E2', '  -0.02%  '),
    @('D3', '1.915.61'),
    @('E3', '  -0.14%  '),
    @('D4', '''1.002'),
    @('E4', '  +0.21%  '),
    @('D5', '''244.64'),
    @('E5', '  -0.68%  '),
    @('E6', '  +0.16%  '),
    @('D7', '''0.4862'),
    @('E7', '  +2.68%  '),
    @('D8', '''0.2903'),
    @('E8', '  +0.58%  '),
    @('D9', '''0.06734'),
    @('E9', '  -1.25%  '),
    @('D10', '''110.96'),
    @('E10', '  +5.33%  '),
    @('D11', '''19.39'),
    @('E11', '  +5.70%  '),
    @('D12', '1.916.96'),
    @('E12', '  -0.04%  '),
    @('D13', '''0.07572'),
    @('E13', '  -1.47%  '),
    @('D14', '''5.368'),
    @('E14', '  +2.01%  '),
    @('D15', '''0.6719'),
    @('E15', '  +0.44%  '),
    @('D16', '''293.06'),
    @('E16', '  +0.33%  '),
    @('D17', '30.580.41'),
    @('E17', '  -0.01%  '),
    @('D18', '''13.03'),
    @('E18', '  +0.85%  '),
    @('E19', '  +0.20%  '),
    @('D20', '''0.000007560'),
    @('E20', '  -0.41%  '),
    @('B21', 'WrappedliquidstakedEther2.0'),
    @('C21', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'),
    @('D21', '2.174.92'),
    @('E21', '  +0.37%  '),
    @('B22', 'Uniswap'),
    @('C22', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'),
    @('D22', '''5.512'),
    @('E22', '  +0.32%  '),
    @('D23', '''1.002'),
    @('D24', '''6.425'),
    @('E24', '  +1.53%  '),
    @('D25', '''9.476'),
    @('E25', '  +1.07%  '),
    @('D26', '''164.76'),
    @('E26', '  -1.98%  '),
    @('D27', '''20.35'),
    @('E27', '  -3.11%  '),
    @('D28', '''2.103'),
    @('E28', '  -0.69%  '),
    @('D29', '''0.1071'),
    @('E29', '  +0.65%  '),
    @('D30', '''1.434'),
    @('E30', '  +2.76%  '),
    @('D31', '''4.139'),
    @('E31', '  -0.56%  '),
    @('D32', '''4.071'),
    @('E32', '  -0.02%  '),
    @('D33', '''0.05008'),
    @('E33', '  -0.48%  '),
    @('D34', '''0.7382'),
    @('E34', '  +0.38%  '),
    @('D35', '''1.138'),
    @('E35', '  -0.59%  '),
    @('D36', '''0.9999'),
    @('E36', '  +0.12%  '),
    @('B37', 'HuobiToken'),
    @('C37', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'),
    @('D37', '''2.716'),
    @('E37', '  -1.12%  '),
    @('B38', 'VeChain'),
    @('C38', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @('D38', '''0.02031'),
    @('E38', '  -1.54%  '),
    @('D39', '''2.685'),
    @('E39', '  -0.08%  '),
    @('D40', '''2.018'),
    @('E40', '  -1.98%  '),
    @('D41', '''109.79'),
    @('E41', '  -1.22%  '),
    @('D42', '''0.4449'),
    @('E42', '  +1.39%  '),
    @('D43', '''0.8639'),
    @('E43', '  -1.50%  '),
    @('D44', '''5.856'),
    @('E44', '  -0.29%  '),
    @('D45', '''69.90'),
    @('E45', '  +4.40%  '),
    @('D46', '''1.002'),
    @('E46', '  +0.22%  '),
    @('D47', '''7.241'),
    @('E47', '  -0.15%  '),
    @('D48', '''9.271'),
    @('E48', '  +0.05%  '),
    @('D49', '''48.11'),
    @('E49', '  +1.15%  '),
    @('D50', '''0.1229'),
    @('E50', '  +0.13%  '),
    @('D51', '''0.2534'),
    @('E51', '  +3.72%  '),
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
